# Generate Report for Handoff
# - Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the rows that just got handed off, and marks their
#   Priority column as "ht" (handoff type) now that a handoff file exists.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 11, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2017-01-03 06:48:26"
}

# --- zh-cn sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2017-01-03 06:48:13"
}

# --- de-de sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2017-01-03 06:48:26"
}
